$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.375.98"
$ws.Range("E2").Value = "  +3.35%  "

$ws.Range("D3").Value = "1.867.88"
$ws.Range("E3").Value = "  +1.72%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "339.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.08%  "

$ws.Range("E6").Value = "  -0.07%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4686"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.72%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3958"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.80%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.30"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.01%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07998"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.61%  "

$ws.Range("E11").Value = "  +2.55%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.84"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.58%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.870.70"
$ws.Range("E13").Value = "  +2.17%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.996"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.81%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.230"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.99%  "

$ws.Range("E16").Value = "  +3.77%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.18%  "

$ws.Range("E18").Value = "  +0.78%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06629"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.19%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.53"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.84%  "

$ws.Range("E21").Value = "  -0.06%  "

$ws.Range("D22").Value = "28.384.21"
$ws.Range("E22").Value = "  +3.37%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.454"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.28%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.268"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.46%  "

$ws.Range("D26").Value = "2.085.96"
$ws.Range("E26").Value = "  +1.52%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "160.60"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.15%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.77"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.31%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.124"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.82%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.484"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.83%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "120.20"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.35%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9650"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.26%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09470"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.19%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.576"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.18%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.342"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.26%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.371"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.90%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06084"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.57%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02242"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.69%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.373"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.83%  "

$ws.Range("E40").Value = "  +2.68%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5932"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.16%  "

$ws.Range("E42").Value = "  +0.02%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1870"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.74%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.34"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.88%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.289"
$ws.Range("D45").Style = "Normal"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5572"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.67%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.14"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.38%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.954"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.85%  "

$ws.Range("E49").Value = "  +3.16%  "

$ws.Range("E50").Value = "  +15.32%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "111.35"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.43%  "
